$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; this shifts existing rows 4-92 down to 5-93
$ws.Rows.Item(4).Insert()

# Fill in the new row 4 with the new price record
$ws.Cells.Item(4, 1).Value = 10
$ws.Cells.Item(4, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(4, 3).Value = "La Araucanía"
$ws.Cells.Item(4, 4).Value = Get-Date -Year 2023 -Month 11 -Day 27 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item(4, 5).Value = 9
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100103
$ws.Cells.Item(4, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(4, 9).Value = 100103003
$ws.Cells.Item(4, 10).Value = "Damasco"
$ws.Cells.Item(4, 11).Value = "Castle Brite"
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 150
$ws.Cells.Item(4, 14).Value = 22000
$ws.Cells.Item(4, 15).Value = 22000
$ws.Cells.Item(4, 16).Value = 22000
$ws.Cells.Item(4, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(4, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(4, 19).Value = 2200
$ws.Cells.Item(4, 20).Value = 10
